# Update the "Approved/Rejected" (column I) and "ReasonToReject" (column J)
# values for the test-case rows (3 through 13) of the Test-Cases sheet.
# Row 2 is left untouched, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 3; $row -le 13; $row++) {
    $ws.Cells.Item($row, 9).Value  = "Rejected"   # column I - Approved/Rejected
    $ws.Cells.Item($row, 10).Value = "Nil"        # column J - ReasonToReject
}

# Reposition the view/selection to match the saved workbook state.
$win = $ws.Application.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("I16").Select()
